$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $text = $cell.Text

    if ([string]::IsNullOrEmpty($text)) {
        continue
    }

    if ($text.Contains(",")) {
        $rawParts = $text.Split(",")
        $parts = @()
        foreach ($p in $rawParts) {
            $parts += $p.Trim()
        }

        $count = $parts.Count
        $lastItem = $parts[$count - 1]
        $rest = $parts[0..($count - 2)]

        $newParts = @($lastItem) + $rest
        $newValue = [string]::Join(", ", $newParts)

        $cell.Value = $newValue
    }
}
